$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B58:B79").Value = ""
